# Add a new "2022-Q1" sheet (fund holdings detail) before the "总计" sheet,
# and add a corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$totalSheetName = "总计"
$templateSheetName = "2021-Q4"

# Locate an existing detail sheet to copy formatting from (2021-Q4, which
# has the same 基金代码/基金名称/... layout as the sheet we are adding).
$templateSheet = $sheets.Item($templateSheetName)

# Insert the new worksheet right before "总计" so the final order is:
# 2021-Q1, 2021-Q3, 2021-Q4, 2022-Q1, 总计
# NOTE: worksheet variables here track sheets by position, so after this
# insertion any previously-fetched reference to the "总计" sheet object
# would now actually point at the freshly inserted sheet instead - always
# re-look-up sheets by name/Item(...) after the collection changes shape.
$newSheet = $sheets.Add($sheets.Item($totalSheetName))
$newSheet.Name = "2022-Q1"

# --- Copy header / column-A cell formatting from the template sheet so the
# new sheet's styles match the other quarterly detail sheets. ---
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Header row ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Data rows ---
# Fund codes (col B, leading zeros) and the numeric-looking values in
# columns D/E/F/G are stored as text in this workbook (e.g. "001499",
# "6.47"), so format those columns as Text before writing them to avoid
# Excel auto-converting them to numbers and dropping leading zeros /
# precision.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$rows = @(
    @(0, "001499", "国投瑞银新增长灵活配置混合A", "6.47", "20.33", "1.23", "0.0796", 1),
    @(1, "011243", "万家惠裕回报6个月持有期混合型证券投资基金A", "4.93", "23.04", "0.78", "0.0385", 9),
    @(2, "007326", "国投瑞银新增长灵活配置混合C", "2.67", "20.33", "1.23", "0.0328", 1),
    @(3, "011244", "万家惠裕回报6个月持有期混合型证券投资基金C", "0.14", "23.04", "0.78", "0.0011", 9)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# --- Update the "总计" (total) summary sheet: insert a new row for
# 2022-Q1 above the existing rows, shifting the other rows' B/C/D data
# down by one. Column A is just a sequential 0-based row index, so it is
# simply renumbered 0..3 for the (now 4) data rows. (Re-fetch the sheet
# by name since the sheet collection shape changed above.) ---
$totalSheet = $sheets.Item($totalSheetName)

# Extend column A's formatting (s=2) down to the new row 5 first.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Shift the existing data rows (B:D) down one row, bottom-up so values
# aren't overwritten before they are read.
$totalSheet.Range("B5").Value = $totalSheet.Range("B4").Value()
$totalSheet.Range("C5").Value = $totalSheet.Range("C4").Value()
$totalSheet.Range("D5").Value = $totalSheet.Range("D4").Value()

$totalSheet.Range("B4").Value = $totalSheet.Range("B3").Value()
$totalSheet.Range("C4").Value = $totalSheet.Range("C3").Value()
$totalSheet.Range("D4").Value = $totalSheet.Range("D3").Value()

$totalSheet.Range("B3").Value = $totalSheet.Range("B2").Value()
$totalSheet.Range("C3").Value = $totalSheet.Range("C2").Value()
$totalSheet.Range("D3").Value = $totalSheet.Range("D2").Value()

# New first data row: 2022-Q1
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.15

# Renumber the column-A row index sequentially for all 4 data rows.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
